$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pairing columns, filled in column-major order (C then D)
$ws.Range("C2").Value = "anette & rich"
$ws.Range("C3").Value = "rick & khali"
$ws.Range("D2").Value = "anette & khali"
$ws.Range("D3").Value = "rick & rich"

# New round headers
$ws.Range("E1").Value = "round 4"
$ws.Range("F1").Value = "round 5"

# Update the selection to match the final state
$ws.Range("B8").Select()
